# Updates the cryptocurrency table (columns B-E, rows 2-51) on Sheet1 to the
# latest scraped price/volume snapshot (GitHub Actions refresh).
#
# Values are written through a throw-away `="literal"` formula and then
# "flattened" in place with Copy + PasteSpecial(xlPasteValues). Writing the
# numeric-looking price strings (e.g. "578.33", "61.895.06") straight into
# .Value would let Excel's type inference convert them into real numbers
# (losing the leading index formatting / introducing float rounding, e.g.
# 578.33 -> 578.33000000000004) or, if forced via NumberFormat="@", stamp a
# new quote-prefixed text style onto the cell. The formula+paste-values round
# trip keeps every touched cell as plain text with its original (default)
# style untouched, matching the source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="62.090.50"'
$ws.Range('E2').Formula = '="  +3.37%  "'
$ws.Range('D3').Formula = '="3.409.25"'
$ws.Range('E3').Formula = '="  +1.74%  "'
$ws.Range('E4').Formula = '="  +0.00%  "'
$ws.Range('D5').Formula = '="578.33"'
$ws.Range('D6').Formula = '="137.76"'
$ws.Range('E6').Formula = '="  +5.18%  "'
$ws.Range('E7').Formula = '="  -0.05%  "'
$ws.Range('E8').Formula = '="  +0.97%  "'
$ws.Range('D9').Formula = '="7.50"'
$ws.Range('E9').Formula = '="  +0.88%  "'
$ws.Range('D10').Formula = '="0.127"'
$ws.Range('E10').Formula = '="  +7.04%  "'
$ws.Range('E11').Formula = '="  +4.33%  "'
$ws.Range('D12').Formula = '="3.988.91"'
$ws.Range('E12').Formula = '="  +1.67%  "'
$ws.Range('E13').Formula = '="  +2.57%  "'
$ws.Range('E14').Formula = '="  +5.50%  "'
$ws.Range('D15').Formula = '="3.404.79"'
$ws.Range('E15').Formula = '="  +1.69%  "'
$ws.Range('D16').Formula = '="25.53"'
$ws.Range('E16').Formula = '="  +3.47%  "'
$ws.Range('D17').Formula = '="62.103.29"'
$ws.Range('E17').Formula = '="  +3.23%  "'
$ws.Range('D18').Formula = '="14.27"'
$ws.Range('E18').Formula = '="  +5.93%  "'
$ws.Range('E19').Formula = '="  +3.26%  "'
$ws.Range('D20').Formula = '="9.50"'
$ws.Range('E20').Formula = '="  +4.44%  "'
$ws.Range('D21').Formula = '="390.04"'
$ws.Range('E21').Formula = '="  +9.97%  "'
$ws.Range('E22').Formula = '="  +2.04%  "'
$ws.Range('D23').Formula = '="3.543.67"'
$ws.Range('E23').Formula = '="  +1.76%  "'
$ws.Range('E24').Formula = '="  +14.88%  "'
$ws.Range('E25').Formula = '="  +0.10%  "'
$ws.Range('D26').Formula = '="71.67"'
$ws.Range('E26').Formula = '="  +3.48%  "'
$ws.Range('D27').Formula = '="7.75"'
$ws.Range('E27').Formula = '="  +2.07%  "'
$ws.Range('E28').Formula = '="  -5.50%  "'
$ws.Range('D29').Formula = '="0.999"'
$ws.Range('E29').Formula = '="  -0.06%  "'
$ws.Range('D30').Formula = '="8.31"'
$ws.Range('E30').Formula = '="  +4.14%  "'
$ws.Range('E31').Formula = '="  +3.45%  "'
$ws.Range('D32').Formula = '="2.18"'
$ws.Range('E32').Formula = '="  +2.48%  "'
$ws.Range('E33').Formula = '="  +0.09%  "'
$ws.Range('D34').Formula = '="3.436.85"'
$ws.Range('E34').Formula = '="  +1.69%  "'
$ws.Range('D35').Formula = '="23.58"'
$ws.Range('E35').Formula = '="  +2.62%  "'
$ws.Range('D36').Formula = '="5.45"'
$ws.Range('E36').Formula = '="  +0.07%  "'
$ws.Range('D37').Formula = '="7.01"'
$ws.Range('E37').Formula = '="  +1.76%  "'
$ws.Range('E38').Formula = '="  +2.81%  "'
$ws.Range('D39').Formula = '="164.12"'
$ws.Range('E39').Formula = '="  +3.81%  "'
$ws.Range('D40').Formula = '="0.0789"'
$ws.Range('E40').Formula = '="  +2.71%  "'
$ws.Range('D41').Formula = '="1.79"'
$ws.Range('E41').Formula = '="  +13.75%  "'
$ws.Range('D42').Formula = '="0.786"'
$ws.Range('E42').Formula = '="  +4.78%  "'
$ws.Range('B43').Formula = '="ONDO"'
$ws.Range('C43').Formula = '="https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"'
$ws.Range('D43').Formula = '="1.24"'
$ws.Range('E43').Formula = '="  +2.74%  "'
$ws.Range('B44').Formula = '="FirstDigitalUSD"'
$ws.Range('C44').Formula = '="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"'
$ws.Range('D44').Formula = '="1.00"'
$ws.Range('E44').Formula = '="  +0.04%  "'
$ws.Range('E45').Formula = '="  +2.21%  "'
$ws.Range('D46').Formula = '="25.15"'
$ws.Range('E46').Formula = '="  +5.82%  "'
$ws.Range('D47').Formula = '="41.70"'
$ws.Range('E47').Formula = '="  +2.34%  "'
$ws.Range('E48').Formula = '="  +2.58%  "'
$ws.Range('D49').Formula = '="23.37"'
$ws.Range('E49').Formula = '="  +3.46%  "'
$ws.Range('D50').Formula = '="2.372.90"'
$ws.Range('E50').Formula = '="  +8.56%  "'
$ws.Range('E51').Formula = '="  +5.73%  "'

# Flatten every formula cell above back down to a literal text value in one
# shot (cells outside the edited set that happen to fall inside this
# bounding range are untouched plain values already, so re-pasting them is a
# no-op).
$changedRange = $ws.Range('B2:E51')
$changedRange.Copy()
$changedRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = 0
